# Fixed update to excel issue
# - Rename "Requested quantity" headers to per-sheet qty names
# - Add a new "PO Forecast" sheet with the Prophet-style forecast output

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the B1 header on "Weekly Quantity" and "Monthly Trend"
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# ---------------------------------------------------------------------------
# 2) Add the new "PO Forecast" sheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- header row ----------------------------------------------------------
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
$headerArr = New-Object 'object[,]' 1, 4
for ($c = 0; $c -lt 4; $c++) { $headerArr[0, $c] = $headers[$c] }
$wsForecast.Range("A1:D1").Value = $headerArr

# Reuse the same bold/bordered/centered header style already present in the
# workbook (copy the formatting straight from an existing header cell so no
# duplicate style entries get created).
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# --- data rows -------------------------------------------------------------
$rows = @(
    @(44983.99999999999, 4, -2.244126452901308, 10.64348337709132),
    @(44997.99999999999, 4, -2.27583753776173, 11.16379785361621),
    @(45004.99999999999, 4, -2.439206241857242, 10.78134173193705),
    @(45011.99999999999, 4, -2.231999088465936, 11.27024074505385),
    @(45039.99999999999, 5, -1.935974279685398, 11.16019196077499),
    @(45088.99999999999, 5, -1.954074941381386, 11.65326961503357),
    @(45095.99999999999, 5, -2.014080279162995, 11.46747338902212),
    @(45102.99999999999, 5, -1.890340310462669, 11.66562955633862),
    @(45193.99999999999, 5, -1.194085217692627, 11.67641288905774),
    @(45200.99999999999, 5, -1.280168822857489, 11.79604676404794),
    @(45207.99999999999, 5, -1.419294311621086, 12.23690086535925),
    @(45221.99999999999, 5, -1.230300399657876, 12.09317282400822),
    @(45235.99999999999, 5, -1.843858059792619, 12.02549685023471),
    @(45452.99999999999, 6, -0.1762936203889576, 13.27233252168687),
    @(45466.99999999999, 6, -0.1412401592257918, 13.15931392279188),
    @(45487.99999999999, 7, 0.1471096151204619, 13.3539271567659),
    @(45508.99999999999, 7, -0.06788119459529098, 13.38674306684116),
    @(45515.99999999999, 7, 0.3529605281263645, 13.33383400890116),
    @(45543.99999999999, 7, 0.02823756136752543, 13.72576800201163),
    @(45578.99999999999, 7, 0.2510134498415148, 13.59747154623267),
    @(45585.99999999999, 7, 0.1480552519907759, 13.52638327431611),
    @(45606.99999999999, 7, 0.4088478424299544, 13.55373399141487),
    @(45634.99999999999, 7, 0.7777151362574042, 13.99904849360332),
    @(45641.99999999999, 7, 0.5288511219649641, 14.1298540070506),
    @(45648.99999999999, 7, 0.4232970835806115, 14.13505892190405),
    @(45655.99999999999, 7, 0.5696060575014007, 13.79968286714684),
    @(45662.99999999999, 7, 0.4228093648581563, 14.24286903027701),
    @(45669.99999999999, 7, 0.8048087318437704, 14.46857923471271),
    @(45676.99999999999, 7, 0.6550632590532101, 13.86668765351263),
    @(45683.99999999999, 7, 0.800932259643614, 14.18364027134656),
    @(45690.99999999999, 7, 0.8016835935513933, 14.21286510297092),
    @(45697.99999999999, 8, 1.125157305484516, 14.02313845064417)
)

$nRows = $rows.Count
$nCols = 4
$dataArr = New-Object 'object[,]' $nRows, $nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $dataArr[$r, $c] = $rows[$r][$c]
    }
}

$endRow = $nRows + 1
$wsForecast.Range("A2:D$endRow").Value = $dataArr

# Column A carries the same date/time format used on the other sheets - copy
# it straight from an existing date cell so the underlying style index (and
# not just the visible format string) matches exactly.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A$endRow").PasteSpecial(-4122)  # xlPasteFormats

# Keep the originally active sheet/selection as it was before this edit.
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()

Write-Output "PO Forecast sheet added with $nRows data rows"
